$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-23"

# Update the header label describing the date range
$ws.Cells.Item(1, 2).Value = "March 2022 (through March 23)"

# Update existing counts that changed
$ws.Cells.Item(5, 2).Value = 7      # B5: 6 -> 7
$ws.Cells.Item(5, 14).Value = 6     # N5: 5 -> 6
$ws.Cells.Item(9, 2).Value = 3      # B9: 4 -> 3
$ws.Cells.Item(25, 14).Value = 3    # N25: 2 -> 3
$ws.Cells.Item(34, 17).Value = 2    # Q34: 1 -> 2
$ws.Cells.Item(41, 2).Value = 4     # B41: 3 -> 4

# Add new counts for previously empty cells
$ws.Cells.Item(14, 8).Value = 1     # H14
$ws.Cells.Item(16, 11).Value = 1    # K16
$ws.Cells.Item(32, 17).Value = 1    # Q32
$ws.Cells.Item(43, 5).Value = 1     # E43
$ws.Cells.Item(44, 17).Value = 1    # Q44
$ws.Cells.Item(46, 14).Value = 1    # N46
$ws.Cells.Item(51, 2).Value = 1     # B51
$ws.Cells.Item(53, 14).Value = 1    # N53
$ws.Cells.Item(56, 8).Value = 1     # H56
$ws.Cells.Item(70, 23).Value = 1    # W70
$ws.Cells.Item(88, 2).Value = 1     # B88
